# Trading update: 2026-02-17 13:18:17
# Appends the newest trade (Trade #21, row 22) to both the "All Trades" and
# "MarketMaking" worksheets, which previously ended at row 21.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 22

    # Trade # (numeric)
    $ws.Cells.Item($row, 1).Value = 21

    # Date / Time — force text so Excel doesn't reinterpret them as
    # date/time serial values.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"

    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = "13:18:13"

    # Strategy / Side
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"

    # Entry Price
    $ws.Cells.Item($row, 6).Value = 0.7

    # Exit Price — trade is still OPEN, so left blank (column G).

    # Status
    $ws.Cells.Item($row, 8).Value = "OPEN"

    # P&L % / P&L $
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0

    # Capital After
    $ws.Cells.Item($row, 11).Value = 99.27951530751794

    # Entry / Exit Slippage (bps)
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0

    # Confidence
    $ws.Cells.Item($row, 14).Value = 0.6

    # Entry Reason
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"

    # Exit Reason — trade still OPEN, so left blank (column P).

    # Duration (min)
    $ws.Cells.Item($row, 17).Value = 0
}
